$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(9285520000000,   118055107.1820449, 2.086868822678343, 3769.001402636258, 3769.001402636258)
    3  = @(9285520000000,   118035702.6683292, 2.083833551784116, 3775.111764408799, 7544.113167045056)
    4  = @(8919010000000,   118032469.7456359, 1.904062509730511, 3968.569170778493, 11512.68233782355)
    5  = @(8919010000000,   118041461.8802993, 1.900953431346346, 3974.757098628297, 15487.43943645185)
    6  = @(10325700000000,  118039409.563591,  1.828997671118578, 4782.768079341661, 20270.20751579351)
    7  = @(10325700000000,  118087151.5187032, 1.83194526512898,  4773.142086607187, 25043.34960240069)
    8  = @(8219280000000,   118066475.638404,  2.079387492158941, 3347.894383041525, 28391.24398544222)
    9  = @(7819060000000,   118100776.9376559, 1.878180692744194, 3525.043001567986, 31916.2869870102)
    10 = @(7953930000000,   118066665.9276808, 1.783929705295721, 3776.389034533244, 35692.67602154345)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
